# Applies the rotation of species-observation rows 3-7 described in the
# commit diff. Each row's core data (Id, Taxonsorteringsordning,
# Rödlistade, TaxonId, Artnamn, Vetenskapligt namn, Auktor, Ost, Nord)
# shifts down by one data row, with the last row (7) wrapping around to
# the first (3). The Start-/Slutdatum (Y/AA) fields only actually change
# value for rows 4 and 6 (rows 3, 5 and 7 keep the same date before and
# after the rotation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 3 (previously held by row 7)
$ws.Range("A3").Value = 80984971
$ws.Range("B3").Value = 89780
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 4217
$ws.Range("F3").Value = "Blodticka"
$ws.Range("G3").Value = "Meruliopsis taxicola"
$ws.Range("H3").Value = "(Pers.:Fr.) Bondartsev"
$ws.Range("Q3").Value = 822333.2752310387
$ws.Range("R3").Value = 7469328.171067456

# New values for row 4 (previously held by row 3)
$ws.Range("A4").Value = 80984949
$ws.Range("B4").Value = 56395
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 823033.1158475638
$ws.Range("R4").Value = 7468549.905579755

# New values for row 5 (previously held by row 4)
$ws.Range("A5").Value = 80984989
$ws.Range("B5").Value = 89952
$ws.Range("D5").Value = "VU"
$ws.Range("E5").Value = 760
$ws.Range("F5").Value = "Doftticka"
$ws.Range("G5").Value = "Haploporus odorus"
$ws.Range("H5").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("Q5").Value = 823861.2863281479
$ws.Range("R5").Value = 7468730.907248246

# New values for row 6 (previously held by row 5)
$ws.Range("A6").Value = 80984988
$ws.Range("B6").Value = 56315
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 102110
$ws.Range("F6").Value = "Fjällvråk"
$ws.Range("G6").Value = "Buteo lagopus"
$ws.Range("H6").Value = "(Pontoppidan, 1763)"
$ws.Range("Q6").Value = 822887.3464848427
$ws.Range("R6").Value = 7468694.206081202

# New values for row 7 (previously held by row 6)
$ws.Range("A7").Value = 80984948
$ws.Range("B7").Value = 56395
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = "Tretåig hackspett"
$ws.Range("G7").Value = "Picoides tridactylus"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("Q7").Value = 823041.2461442947
$ws.Range("R7").Value = 7468541.086727832

# Startdatum/Slutdatum text values only actually swap for rows 4 and 6;
# the cells hold plain text (not real dates), so force Text format
# before assignment to avoid Excel auto-converting the string into a
# date serial number, then restore the Normal style so the cell's
# formatting stays exactly as it was originally (General/default).
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2019-06-19"
$ws.Range("Y4").Style = "Normal"

$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2019-06-19"
$ws.Range("AA4").Style = "Normal"

$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2019-06-18"
$ws.Range("Y6").Style = "Normal"

$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2019-06-18"
$ws.Range("AA6").Style = "Normal"
